# Refresh the crypto price/volume snapshot (GitHub Actions data pull).
# Column D ("Price") and most of column E ("Volume(1h)") are stored as
# plain text (the source feed already formats them, e.g. "70.641.38" or
# "  +2.07%  "), and rows 19/20 swapped rank (Uniswap <-> Polkadot) while
# row 51 changed from FLOKI to Cosmos.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.641.38"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "3.808.66"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "668.43"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +7.03%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "169.59"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("D7").Value = "3.805.21"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("E10").Value = "  +0.45%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.463"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("E12").Value = "  +4.48%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000245"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.70%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "35.94"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "4.457.27"
$ws.Range("D16").Value = "3.814.00"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "70.658.01"
$ws.Range("E17").Value = "  +1.98%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "17.73"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.70"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +21.65%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.21"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("E21").Value = "  +0.61%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "475.07"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.36%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.716"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.54%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "83.14"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -2.32%  "
$ws.Range("E26").Value = "  +1.71%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.37"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.67%  "
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "3.964.85"
$ws.Range("E30").Value = "  +0.89%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.87"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +7.70%  "
$ws.Range("E32").Value = "  +2.60%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "7.40"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.80%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "29.71"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.21%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.175"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +5.21%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "9.16"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "3.769.86"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("E39").Value = "  +0.63%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.44"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.90%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.02"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.74%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  +9.44%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "45.75"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +5.91%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "158.66"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.55%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "47.96"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.67%  "
$ws.Range("E49").Value = "  +4.62%  "
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "8.53"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.18%  "
